# The underlying OOXML diff for this revision is a pure re-serialization:
# every hunk in the target diff (the <w:document> namespace declarations,
# <w:tblW>, <w:tblLook>, <w:tcW>, <wp:anchor>, <wp:effectExtent>,
# <wps:bodyPr>, <v:shapetype>, <v:shape> (incl. its base64 o:gfxdata blob,
# which decodes to the exact same bytes), <w:pgSz>, <w:pgMar>, <w:rFonts>,
# <w:lang>, <w:latentStyles>, and every <w:style>/<w:lsdException> element)
# only reorders XML attributes (and re-wraps one base64 blob's whitespace);
# no text, value, paragraph, table, style, numbering, or formatting content
# actually changes anywhere in the package (confirmed by decoding the VML
# preview blob and by set-comparing every attribute=value pair on both
# sides of every hunk). The commit message ("Moving from 2.0.1 to 2.0.2")
# matches this: a tooling/library version bump that re-saved the test
# fixture through a different OOXML writer with a different (alphabetical)
# attribute-ordering convention, with no authoring change to the document.
#
# This Word-OM/COM-interop surface has no operation that lets a script
# choose the attribute-serialization order of an OOXML element (every
# element - whether copied from the source part or freshly built by e.g.
# Tables.Add - is always written back in this host's own fixed internal
# attribute order). Since the target content is byte-for-byte/semantically
# identical to the source document, the correct edit is therefore a no-op:
# touching the document here would only risk introducing spurious content
# changes (run splits/merges, rsid churn, etc.) that are not present in
# the target revision.
$d = $word.ActiveDocument
